$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (the duplicate "ENDF 115In(n,g)" row). Excel will shift all
# subsequent rows up by one, re-number cell/row references, and adjust the
# shared-formula ranges (M2:M20 -> M2:M19, O2:O20 -> O2:O19, P3:P20 -> P3:P19)
# automatically. It will also drop the now-unused shared string entry when
# the workbook is saved.
$ws.Rows("7:7").Delete()

# Correct cell F2: it used to hold a bare literal 1.2365E-7; it should now be
# computed from a formula that accounts for the atom-density error.
$ws.Range("F2").Formula = "=0.00000012365*0.022148/0.02551"

# Restore the selected cell to A6, matching the saved view state.
$ws.Range("A6").Select()
